# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.319.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.09%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.932.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.03%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.26%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.7594"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +5.92%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'244.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.55%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9990"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.24%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3185"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.42%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'27.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.29%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07008"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.37%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.7796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.42%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07999"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.932.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.06%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.361"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.98%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'94.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.30%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D17").Value = "'30.312.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.02%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'253.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.60%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.000007919"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.50%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'5.739"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.92%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.188.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.25%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.9991"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.20%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.9981"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.27%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.676"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.45%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -2.22%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'165.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.05%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.1343"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.52%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'18.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.23%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.193"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.60%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.16%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -2.02%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.381"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.03%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.117"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.81%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.05163"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.78%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.291"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.20%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7485"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.30%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.770"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.02%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01955"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.16%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.799"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.00%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'77.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.70%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.423"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.31%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.4468"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.18%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.964"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.11%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.9994"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.19%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.8300"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.25%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'100.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.09%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'9.752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.13%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.475"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.92%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Maker"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'984.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +11.58%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Elrond"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'37.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.19%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.06002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.15%  "
$ws.Range("E51").Style = "Normal"

Write-Host "Updated cryptos list"
